$d = $word.ActiveDocument

# Remove ", data_di_nascita" from the Collezionista tuple definition
$d.Content.Find.Execute(
    "Collezionista(ID, nickname, email, nome, cognome, data_di_nascita)",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Collezionista(ID, nickname, email, nome, cognome)",
    2
)

# Remove "data_di_nascita, " from the Artista tuple definition
$d.Content.Find.Execute(
    "Artista(ID, nome_dArte, nome, cognome, data_di_nascita, gruppo)",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Artista(ID, nome_dArte, nome, cognome, gruppo)",
    2
)
